$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add row 3 with A3 = AVERAGE(A1:A2)
$ws.Range("A3").Formula = "=AVERAGE(A1:A2)"

# Update selection to A4
$ws.Range("A4").Select()
